$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Produtos")
$ws.Range("I2").ClearContents()
$ws.Range("I2").Select()
